$p = $ppt.ActivePresentation
$s5 = $p.Slides.Item(5)
$rect = $s5.Shapes.Item("Rectangle 84")
$tr = $rect.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf('"OSMTILE"')
$run = $tr.Characters($idx + 1, 9)
$lq = [char]0x201C
$rq = [char]0x201D
$run.Text = '"OSMTILE' + $lq + ' checked=' + $lq + 'checked' + $rq
$rect.Height = 327.16405511811024
Write-Output ("Height after fix=" + $rect.Height)
